# cryptos.xlsx -- scheduled data refresh (GitHub Actions bot).
# Row order shuffles slightly because the source list is re-ranked by
# volume each run (Cosmos/EthereumClassic swap @28-29, TrustWalletToken/
# FraxShare swap @42-43, and BabyDogeCoin enters @48, pushing Cronos/
# EnergySwap/Mantle down one row and dropping USDD off the bottom).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = "'26.694.01"
$ws.Cells.Item(2, 5).Value = '  +1.42%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = "'1.634.98"
$ws.Cells.Item(3, 5).Value = '  +1.36%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = '  +0.06%  '

# Row 5: BNB
$ws.Cells.Item(5, 4).Value = "'213.41"
$ws.Cells.Item(5, 5).Value = '  +0.13%  '

# Row 6: USDC
$ws.Cells.Item(6, 5).Value = '  +0.01%  '

# Row 7: XRP
$ws.Cells.Item(7, 5).Value = '  +1.03%  '

# Row 8: Cardano
$ws.Cells.Item(8, 5).Value = '  +0.40%  '

# Row 9: Dogecoin
$ws.Cells.Item(9, 5).Value = '  +0.66%  '

# Row 10: Solana
$ws.Cells.Item(10, 5).Value = '  +2.67%  '

# Row 11: TRON
$ws.Cells.Item(11, 5).Value = '  +2.42%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Cells.Item(12, 4).Value = "'1.861.92"
$ws.Cells.Item(12, 5).Value = '  +1.35%  '

# Row 13: WrappedEther
$ws.Cells.Item(13, 4).Value = "'1.632.37"
$ws.Cells.Item(13, 5).Value = '  +1.32%  '

# Row 14: Polkadot
$ws.Cells.Item(14, 5).Value = '  +0.06%  '

# Row 15: Polygon
$ws.Cells.Item(15, 5).Value = '  +1.40%  '

# Row 16: WrappedBTC
$ws.Cells.Item(16, 4).Value = "'26.665.24"
$ws.Cells.Item(16, 5).Value = '  +1.31%  '

# Row 17: Litecoin
$ws.Cells.Item(17, 4).Value = "'63.14"
$ws.Cells.Item(17, 5).Value = '  +2.01%  '

# Row 18: ShibaInu
$ws.Cells.Item(18, 5).Value = '  +0.54%  '

# Row 19: Dai
$ws.Cells.Item(19, 5).Value = '  +0.07%  '

# Row 20: BitcoinCash
$ws.Cells.Item(20, 4).Value = "'208.39"
$ws.Cells.Item(20, 5).Value = '  +2.50%  '

# Row 21: Uniswap
$ws.Cells.Item(21, 5).Value = '  +0.42%  '

# Row 22: Avalanche
$ws.Cells.Item(22, 4).Value = "'9.37"
$ws.Cells.Item(22, 5).Value = '  +0.48%  '

# Row 23: Chainlink
$ws.Cells.Item(23, 4).Value = "'6.09"
$ws.Cells.Item(23, 5).Value = '  +0.90%  '

# Row 24: Toncoin
$ws.Cells.Item(24, 4).Value = "'1.90"
$ws.Cells.Item(24, 5).Value = '  -0.53%  '

# Row 25: Monero
$ws.Cells.Item(25, 4).Value = "'145.78"
$ws.Cells.Item(25, 5).Value = '  +1.01%  '

# Row 26: BinanceUSD
$ws.Cells.Item(26, 5).Value = '  +0.06%  '

# Row 27: Stellar
$ws.Cells.Item(27, 5).Value = '  -1.59%  '

# Row 28: Cosmos
$ws.Cells.Item(28, 2).Value = 'Cosmos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(28, 4).Value = "'6.68"
$ws.Cells.Item(28, 5).Value = '  +1.60%  '

# Row 29: EthereumClassic
$ws.Cells.Item(29, 2).Value = 'EthereumClassic'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(29, 4).Value = "'15.35"
$ws.Cells.Item(29, 5).Value = '  +0.59%  '

# Row 30: Hedera
$ws.Cells.Item(30, 5).Value = '  +5.58%  '

# Row 31: PancakeSwap
$ws.Cells.Item(31, 5).Value = '  -0.15%  '

# Row 32: Filecoin
$ws.Cells.Item(32, 4).Value = "'3.23"
$ws.Cells.Item(32, 5).Value = '  +0.68%  '

# Row 33: InternetComputer(DFINITY)
$ws.Cells.Item(33, 4).Value = "'2.95"
$ws.Cells.Item(33, 5).Value = '  -0.42%  '

# Row 34: LidoDAOToken
$ws.Cells.Item(34, 5).Value = '  +1.04%  '

# Row 35: HuobiToken
$ws.Cells.Item(35, 5).Value = '  -0.35%  '

# Row 36: Maker
$ws.Cells.Item(36, 4).Value = "'1.164.61"
$ws.Cells.Item(36, 5).Value = '  +0.24%  '

# Row 37: VeChain
$ws.Cells.Item(37, 5).Value = '  +0.82%  '

# Row 38: ARBITRUM
$ws.Cells.Item(38, 4).Value = "'0.812"
$ws.Cells.Item(38, 5).Value = '  +1.85%  '

# Row 39: PaxDollar
$ws.Cells.Item(39, 5).Value = '  +0.04%  '

# Row 40: MXToken
$ws.Cells.Item(40, 5).Value = '  +0.01%  '

# Row 41: ImmutableX
$ws.Cells.Item(41, 5).Value = '  -0.10%  '

# Row 42: TrustWalletToken
$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 4).Value = "'0.795"
$ws.Cells.Item(42, 5).Value = '  +0.91%  '

# Row 43: FraxShare
$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43, 4).Value = "'5.40"
$ws.Cells.Item(43, 5).Value = '  +2.65%  '

# Row 44: RocketPoolETH
$ws.Cells.Item(44, 4).Value = "'1.772.82"
$ws.Cells.Item(44, 5).Value = '  +1.37%  '

# Row 45: Quant
$ws.Cells.Item(45, 4).Value = "'92.38"
$ws.Cells.Item(45, 5).Value = '  +0.53%  '

# Row 46: RenderToken
$ws.Cells.Item(46, 4).Value = "'1.55"
$ws.Cells.Item(46, 5).Value = '  +0.52%  '

# Row 47: Aave
$ws.Cells.Item(47, 4).Value = "'54.67"
$ws.Cells.Item(47, 5).Value = '  +0.43%  '

# Row 48: BabyDogeCoin
$ws.Cells.Item(48, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(48, 4).Value = "'0.0₆0101"
$ws.Cells.Item(48, 5).Value = '  +4.22%  '

# Row 49: Cronos
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).Value = "'0.0512"
$ws.Cells.Item(49, 5).Value = '  +0.83%  '

# Row 50: EnergySwap
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = "'7.58"
$ws.Cells.Item(50, 5).Value = '  +4.47%  '

# Row 51: Mantle
$ws.Cells.Item(51, 2).Value = 'Mantle'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(51, 4).Value = "'0.410"
$ws.Cells.Item(51, 5).Value = '  +0.77%  '
